$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "山子高科"
$ws.Range("C2").Value = "山子高科"
$ws.Range("A3").Value = "山子高科"
$ws.Range("B3").Value = "北方稀土"
$ws.Range("C3").Value = "闻泰科技"
$ws.Range("A4").Value = "包钢股份"
$ws.Range("B4").Value = "包钢股份"
$ws.Range("C4").Value = "通富微电"
$ws.Range("A5").Value = "安泰科技"
$ws.Range("C5").Value = "安泰科技"
$ws.Range("A6").Value = "至纯科技"
$ws.Range("C6").Value = "至纯科技"
$ws.Range("A7").Value = "合锻智能"
$ws.Range("B7").Value = "合锻智能"
$ws.Range("C7").Value = "合锻智能"
$ws.Range("A8").Value = "新莱应材"
$ws.Range("B8").Value = "长城军工"
$ws.Range("C8").Value = "北方稀土"
$ws.Range("A9").Value = "闻泰科技"
$ws.Range("B9").Value = "白银有色"
$ws.Range("C9").Value = "凯美特气"
$ws.Range("A10").Value = "南 京 港"
$ws.Range("B10").Value = "凯美特气"
$ws.Range("C10").Value = "长城军工"
$ws.Range("B11").Value = "安泰科技"
$ws.Range("C11").Value = "新莱应材"
$ws.Range("A12").Value = "凯美特气"
$ws.Range("B12").Value = "新莱应材"
$ws.Range("C12").Value = "盛和资源"
$ws.Range("A13").Value = "通富微电"
$ws.Range("B13").Value = "隆基绿能"
$ws.Range("C13").Value = "包钢股份"
$ws.Range("A14").Value = "白银有色"
$ws.Range("B14").Value = "楚江新材"
$ws.Range("C14").Value = "白银有色"
$ws.Range("A15").Value = "隆基绿能"
$ws.Range("B15").Value = "通富微电"
$ws.Range("C15").Value = "合肥城建"
$ws.Range("A16").Value = "黄河旋风"
$ws.Range("B16").Value = "中油资本"
$ws.Range("C16").Value = "楚江新材"
$ws.Range("A17").Value = "楚江新材"
$ws.Range("B17").Value = "中国稀土"
$ws.Range("C17").Value = "中电鑫龙"
$ws.Range("A18").Value = "中国稀土"
$ws.Range("B18").Value = "南 京 港"
$ws.Range("C18").Value = "南京港"
$ws.Range("A19").Value = "上海电气"
$ws.Range("B19").Value = "上海电气"
$ws.Range("C19").Value = "上海电气"
$ws.Range("A20").Value = "合肥城建"
$ws.Range("B20").Value = "贵州茅台"
$ws.Range("C20").Value = "黄河旋风"
$ws.Range("A21").Value = "中电鑫龙"
$ws.Range("B21").Value = "黄河旋风"
$ws.Range("C21").Value = "领益智造"
